# Update "想去人数" (number of people interested) values in column F
# for the sheets "展览" and "全部类型" which contain identical data.

$wb = $excel.ActiveWorkbook

# Mapping of row number -> new value for column F
$updates = @{
    2  = 1942
    7  = 1615
    9  = 639
    11 = 101
    12 = 21
    17 = 108
    19 = 3750
    22 = 431
    23 = 345
    24 = 701
    25 = 405
    26 = 352
    28 = 1538
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
